$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("M2").Value = 1.543058333333333
$ws.Range("N2").Value = 4.629175
$ws.Range("O2").Value = 0.3927712126299722
$ws.Range("P2").Value = 0.4267890789273644
$ws.Range("Q2").Value = 2.4797934974
$ws.Range("R2").Value = 22.3181414766
$ws.Range("S2").Value = 0.3927712126299722
$ws.Range("T2").Value = 0.4267890789273644

$ws.Range("M3").Value = 1.446170333333334
$ws.Range("O3").Value = 0.3681092692495906
$ws.Range("P3").Value = 0.399991167671613
$ws.Range("Q3").Value = 2.324088280568001
$ws.Range("S3").Value = 0.3681092692495906
$ws.Range("T3").Value = 0.399991167671613

$ws.Range("M4").Value = 0.9394155
$ws.Range("N4").Value = 1.878831
$ws.Range("O4").Value = 0.2391195181204372
$ws.Range("P4").Value = 0.1732197534010226
$ws.Range("Q4").Value = 1.509700831092
$ws.Range("R4").Value = 9.058204986551999
$ws.Range("S4").Value = 0.2391195181204372
$ws.Range("T4").Value = 0.1732197534010226
